$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = -0.7142215300419087

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0

$ws.Range("B4").Value = 0
$ws.Range("C4").Value = -0.7495068380493121
$ws.Range("D4").Value = 0

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = -0.6600173190601064
$ws.Range("D5").Value = 0.7349814711053394

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0.8285788005719147

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0

$ws.Range("B8").Value = 0
$ws.Range("C8").Value = -0.6925154799362028
$ws.Range("D8").Value = 0

$ws.Range("B9").Value = 0
$ws.Range("C9").Value = -0.8107899352968411
$ws.Range("D9").Value = 0.7074476980345612
